$wb = $excel.ActiveWorkbook

# --- Sheet: ip_address_list ---
$ws1 = $wb.Worksheets.Item("ip_address_list")

# B1: fix an IP typo
$ws1.Range("B1").Value = "192.168.10.255"

# B2: new IP, E2: number -> boolean (keep value 1/TRUE)
$ws1.Range("B2").Value = "192.168.18.241"
$ws1.Range("E2").Value = $true

# B3: append trailing "g" (typo), E3: number -> boolean (keep value 0/FALSE)
$ws1.Range("B3").Value = "192.168.000.000g"
$ws1.Range("E3").Value = $false

# --- Sheet: ip_adress_fav_list ---
$ws2 = $wb.Worksheets.Item("ip_adress_fav_list")

# Swap row 1 and row 2 entirely (keeps original cell types/formatting intact,
# e.g. A1/A2 stay text even though "514" looks numeric) via a scratch row.
$ws2.Range("A1:E1").Copy($ws2.Range("A10:E10"))
$ws2.Range("A2:E2").Copy($ws2.Range("A1:E1"))
$ws2.Range("A10:E10").Copy($ws2.Range("A2:E2"))
$ws2.Range("A10:E10").Clear()

# After the swap, row2's IP (originally row1's IP) gets replaced with a new one.
$ws2.Range("B2").Value = "192.168.18.241"
